# Prefix each protocol step's Column A label with its sheet (protocol) name.
# e.g. on sheet "discount2", "Step4 Takeaway" -> "discount2 Step4 Takeaway"
# Row 1 holds the header "Name" and must remain untouched.

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "price1", "price2",
    "discount1", "discount2",
    "free1", "free2",
    "nomoney1", "nomoney2",
    "noppv1", "noppv2",
    "card1", "card2",
    "nosex1", "nosex2",
    "offtopic1", "offtopic2",
    "real1", "real2",
    "voice1", "voice2",
    "customyes1", "customyes2",
    "customno1", "customno2",
    "done1", "done2",
    "cumcontrol", "dickpic", "boosters"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ur = $ws.UsedRange
    $lastRow = $ur.Row + $ur.Rows.Count - 1

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Value2
        if ($current -ne $null -and $current -ne "") {
            $cell.Value = "$name $current"
        }
    }
}
